# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The "K" column (G) values were recalculated/regenerated; write the new values
# for rows 2-76 into column G (header "K" at G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(
    1,2,1,1,3,0,2,1,0,2,0,1,1,2,2,0,0,2,1,0,
    1,1,3,3,1,2,2,0,0,2,1,1,2,1,0,3,2,1,0,2,
    2,1,1,3,1,1,1,0,1,0,3,2,0,1,2,0,2,0,1,0,
    3,2,1,1,2,1,2,1,2,1,1,3,0,2,0
)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
